# Case_5_201 (380 kV) results: refresh res_bus/vm_pu.xlsx values for the
# new case run. Column B (bus 0, the slack/external-grid bus) moves from
# the old 1.05 pu setpoint to 1.02 pu, and every downstream bus voltage
# (columns C-F, I-M) is updated to the newly computed power-flow results.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "B2" = 1.02
    "C2" = 1.067679594307383
    "D2" = 1.066407713249685
    "E2" = 1.071452302287769
    "F2" = 1.079542135599915
    "I2" = 1.046038474011035
    "J2" = 1.072623263039359
    "K2" = 1.069118343834025
    "L2" = 1.07414943788308
    "M2" = 1.082217914197944
    "B3" = 1.02
    "C3" = 1.069114138464886
    "D3" = 1.06748589851349
    "E3" = 1.072692109069219
    "F3" = 1.080775892446987
    "I3" = 1.046348868933508
    "J3" = 1.073711991790586
    "K3" = 1.070011460488645
    "L3" = 1.075204755659586
    "M3" = 1.08326875024748
    "B4" = 1.02
    "C4" = 1.07004180243142
    "D4" = 1.06818283535275
    "E4" = 1.073494016840617
    "F4" = 1.081573801282443
    "I4" = 1.046548065051609
    "J4" = 1.074415421533423
    "K4" = 1.070588047754984
    "L4" = 1.075886719912771
    "M4" = 1.083947726777279
    "B5" = 1.02
    "C5" = 1.070431657631805
    "D5" = 1.06847565750326
    "E5" = 1.073831062954426
    "F5" = 1.081909145960604
    "I5" = 1.046631413012683
    "J5" = 1.074710895026525
    "K5" = 1.070830131727059
    "L5" = 1.076173205838833
    "M5" = 1.084232935555082
    "B6" = 1.02
    "C6" = 1.070497108268382
    "D6" = 1.068524813729912
    "E6" = 1.073887650118832
    "F6" = 1.081965446253594
    "I6" = 1.046645384420545
    "J6" = 1.074760491852502
    "K6" = 1.070870760388987
    "L6" = 1.076221295747096
    "M6" = 1.084280809783682
    "B7" = 1.02
    "C7" = 1.07004701221936
    "D7" = 1.068186748721873
    "E7" = 1.073498520761702
    "F7" = 1.08157828254871
    "I7" = 1.046549180298447
    "J7" = 1.07441937063674
    "K7" = 1.07059128372096
    "L7" = 1.075890548780579
    "M7" = 1.083951538662844
    "B8" = 1.02
    "C8" = 1.068164528511769
    "D8" = 1.066772241854072
    "E8" = 1.071871371027985
    "F8" = 1.079959177103673
    "I8" = 1.046143715684216
    "J8" = 1.07299142343431
    "K8" = 1.069420451353884
    "L8" = 1.074506274913979
    "M8" = 1.082573254714964
    "B9" = 1.02
    "C9" = 1.064842687390796
    "D9" = 1.064274062646977
    "E9" = 1.069001456812278
    "F9" = 1.077102798208031
    "I9" = 1.045416553610883
    "J9" = 1.070467011085585
    "K9" = 1.067347083910478
    "L9" = 1.072060008845984
    "M9" = 1.08013687617183
    "B10" = 1.02
    "C10" = 1.062624699915224
    "D10" = 1.062604657164433
    "E10" = 1.067086188287852
    "F10" = 1.0751961303888
    "I10" = 1.04492319059676
    "J10" = 1.068778366064114
    "K10" = 1.065957817945278
    "L10" = 1.070424276395139
    "M10" = 1.078507292869472
    "B11" = 1.02
    "C11" = 1.061663408848762
    "D11" = 1.061880811749488
    "E11" = 1.066256338498263
    "F11" = 1.074369906085119
    "I11" = 1.044707507780898
    "J11" = 1.06804576642748
    "K11" = 1.06535454990176
    "L11" = 1.069714785383795
    "M11" = 1.077800360891305
    "B12" = 1.02
    "C12" = 1.061306203454833
    "D12" = 1.061611792551954
    "E12" = 1.065948012381001
    "F12" = 1.074062911940763
    "I12" = 1.044627083802169
    "J12" = 1.0677734312294
    "K12" = 1.065130209748009
    "L12" = 1.069451063685283
    "M12" = 1.077537574279896
    "B13" = 1.02
    "C13" = 1.061382831610344
    "D13" = 1.061669504958972
    "E13" = 1.066014153197381
    "F13" = 1.074128767707264
    "I13" = 1.044644349057355
    "J13" = 1.067831857850194
    "K13" = 1.065178343267756
    "L13" = 1.069507641321918
    "M13" = 1.077593952034611
    "B14" = 1.02
    "C14" = 1.06163388500994
    "D14" = 1.061858577633595
    "E14" = 1.066230853880432
    "F14" = 1.074344531867579
    "I14" = 1.044700866236536
    "J14" = 1.06802325952547
    "K14" = 1.065336011187399
    "L14" = 1.069692989864042
    "M14" = 1.077778642980711
    "B15" = 1.02
    "C15" = 1.061788548693428
    "D15" = 1.061975051470383
    "E15" = 1.066364359212259
    "F15" = 1.074477458233664
    "I15" = 1.044735647248933
    "J15" = 1.068141159796479
    "K15" = 1.065433121110703
    "L15" = 1.069807164574882
    "M15" = 1.077892410480759
    "B16" = 1.02
    "C16" = 1.062688478445375
    "D16" = 1.062652675435859
    "E16" = 1.067141251269107
    "F16" = 1.07525095069068
    "I16" = 1.044937461383469
    "J16" = 1.068826956345889
    "K16" = 1.065997818633489
    "L16" = 1.07047133721338
    "M16" = 1.078554181680112
    "B17" = 1.02
    "C17" = 1.063252738882951
    "D17" = 1.063077465879555
    "E17" = 1.067628431556768
    "F17" = 1.075735972173208
    "I17" = 1.045063503408981
    "J17" = 1.069256759402474
    "K17" = 1.066351579332153
    "L17" = 1.070887629158737
    "M17" = 1.07896893961104
    "B18" = 1.02
    "C18" = 1.063581777172037
    "D18" = 1.063325144449753
    "E18" = 1.067912545105212
    "F18" = 1.076018816997506
    "I18" = 1.045136823556443
    "J18" = 1.069507320941725
    "K18" = 1.0665577575025
    "L18" = 1.071130328894869
    "M18" = 1.079210734580221
    "B19" = 1.02
    "C19" = 1.063693956432598
    "D19" = 1.063409580509604
    "E19" = 1.068009412041063
    "F19" = 1.076115249710844
    "I19" = 1.045161790280513
    "J19" = 1.069592733119342
    "K19" = 1.066628031129199
    "L19" = 1.071213063614064
    "M19" = 1.079293159097327
    "B20" = 1.02
    "C20" = 1.063192207941635
    "D20" = 1.063031899685453
    "E20" = 1.067576166972855
    "F20" = 1.075683940193396
    "I20" = 1.045050000785792
    "J20" = 1.069210659627973
    "K20" = 1.066313641191712
    "L20" = 1.070842977031089
    "M20" = 1.078924453102397
    "B21" = 1.02
    "C21" = 1.061559959926592
    "D21" = 1.06180290463212
    "E21" = 1.066167043252221
    "F21" = 1.074280997426545
    "I21" = 1.044684231897627
    "J21" = 1.067966902471658
    "K21" = 1.065289589104668
    "L21" = 1.069638414478438
    "M21" = 1.077724261666397
    "B22" = 1.02
    "C22" = 1.060532893239586
    "D22" = 1.061029312654418
    "E22" = 1.065280589309959
    "F22" = 1.073398344745313
    "I22" = 1.044452465562613
    "J22" = 1.067183656655331
    "K22" = 1.06464422372204
    "L22" = 1.068879985169137
    "M22" = 1.076968491307627
    "B23" = 1.02
    "C23" = 1.061077439181403
    "D23" = 1.061439492341994
    "E23" = 1.065750562321907
    "F23" = 1.073866310645375
    "I23" = 1.044575499659265
    "J23" = 1.067598989546025
    "K23" = 1.06498648765551
    "L23" = 1.069282145828865
    "M23" = 1.077369250838985
    "B24" = 1.02
    "C24" = 1.063219559531654
    "D24" = 1.06305248937909
    "E24" = 1.067599783244402
    "F24" = 1.075707451390623
    "I24" = 1.045056102651841
    "J24" = 1.069231490549623
    "K24" = 1.066330784312363
    "L24" = 1.070863153764761
    "M24" = 1.078944555030889
    "B25" = 1.02
    "C25" = 1.065702045544786
    "D25" = 1.064920586572346
    "E25" = 1.069743735003671
    "F25" = 1.077841652257167
    "I25" = 1.045606051622895
    "J25" = 1.071120622810169
    "K25" = 1.067884325135804
    "L25" = 1.072693275409668
    "M25" = 1.080767665029961
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output "Updated $($updates.Count) cells"
